# Update the "time_taken" (column F) timestamps on the existing "data" sheet
# to reflect the re-run query time (10:51:30.9xxxxx -> 14:34:24.5xxxxx).
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$timeTaken = @{
    2  = "2021-10-05 14:34:24.540447"
    3  = "2021-10-05 14:34:24.540455"
    4  = "2021-10-05 14:34:24.540458"
    5  = "2021-10-05 14:34:24.540461"
    6  = "2021-10-05 14:34:24.540463"
    7  = "2021-10-05 14:34:24.540466"
    8  = "2021-10-05 14:34:24.540469"
    9  = "2021-10-05 14:34:24.540471"
    10 = "2021-10-05 14:34:24.540474"
    11 = "2021-10-05 14:34:24.540477"
    12 = "2021-10-05 14:34:24.540479"
    13 = "2021-10-05 14:34:24.540482"
    14 = "2021-10-05 14:34:24.540485"
    15 = "2021-10-05 14:34:24.540487"
    16 = "2021-10-05 14:34:24.540490"
    17 = "2021-10-05 14:34:24.540492"
    18 = "2021-10-05 14:34:24.540495"
    19 = "2021-10-05 14:34:24.540498"
    20 = "2021-10-05 14:34:24.540500"
    21 = "2021-10-05 14:34:24.540503"
    22 = "2021-10-05 14:34:24.540505"
    23 = "2021-10-05 14:34:24.540508"
    24 = "2021-10-05 14:34:24.540510"
    25 = "2021-10-05 14:34:24.540513"
    26 = "2021-10-05 14:34:24.540515"
    27 = "2021-10-05 14:34:24.540518"
    28 = "2021-10-05 14:34:24.540521"
    29 = "2021-10-05 14:34:24.540523"
    30 = "2021-10-05 14:34:24.540526"
    31 = "2021-10-05 14:34:24.540528"
    32 = "2021-10-05 14:34:24.540531"
    33 = "2021-10-05 14:34:24.540533"
    34 = "2021-10-05 14:34:24.540536"
    35 = "2021-10-05 14:34:24.540539"
    36 = "2021-10-05 14:34:24.540541"
    37 = "2021-10-05 14:34:24.540544"
    38 = "2021-10-05 14:34:24.540546"
    39 = "2021-10-05 14:34:24.540549"
    40 = "2021-10-05 14:34:24.540551"
    41 = "2021-10-05 14:34:24.540554"
    42 = "2021-10-05 14:34:24.540556"
    43 = "2021-10-05 14:34:24.540559"
    44 = "2021-10-05 14:34:24.540562"
    45 = "2021-10-05 14:34:24.540564"
    46 = "2021-10-05 14:34:24.540567"
    47 = "2021-10-05 14:34:24.540569"
    48 = "2021-10-05 14:34:24.540572"
    49 = "2021-10-05 14:34:24.540574"
    50 = "2021-10-05 14:34:24.540577"
    51 = "2021-10-05 14:34:24.540579"
    52 = "2021-10-05 14:34:24.540582"
    53 = "2021-10-05 14:34:24.540584"
    54 = "2021-10-05 14:34:24.540587"
    55 = "2021-10-05 14:34:24.540590"
    56 = "2021-10-05 14:34:24.540592"
    57 = "2021-10-05 14:34:24.540595"
    58 = "2021-10-05 14:34:24.540597"
    59 = "2021-10-05 14:34:24.540600"
    60 = "2021-10-05 14:34:24.540602"
    61 = "2021-10-05 14:34:24.540605"
    62 = "2021-10-05 14:34:24.540607"
}

foreach ($row in $timeTaken.Keys) {
    $dataSheet.Cells.Item($row, 6).Value = $timeTaken[$row]
}

# Add the new "metadata" sheet (placed after "data") capturing the panel
# query metadata for this run.
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Match header styling used on the "data" sheet's header row
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Liver Failure_Paediatric"
$ws.Range("C2").Value = 3400

# data_version must stay text ("1.8"), not be coerced to a number
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.8"
$dataSheet.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("E2").Value = "2021-09-27T07:40:54.623315Z"
$ws.Range("F2").Value = "2021-10-05 14:34:24.536740"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3400/?format=json"

# Match the numeric-cell styling used on the "data" sheet's A column
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

Write-Output "done"
